# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a new "Label" column (H) marking Control (0) vs MDD (1) rows,
# and refreshes a handful of recomputed Prediction/Error/Cross-Entropy values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" header in H1, styled like the other headers (copy format from G1) ---
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- New "Label" data column: 0 for Control rows, 1 for MDD rows (both blocks) ---
$labels = @(0, 0, 0, 0, 0, 1, 1, 1, 1, 1)
for ($i = 0; $i -lt 10; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $row = 12 + $i
    $ws.Cells.Item($row, 8).Value = $labels[$i]
}

# --- Refreshed Prediction / Error / Cross Entropy Loss values (re-fit results) ---
$ws.Range("D5").Value = 0.48603924639782
$ws.Range("E5").Value = 0.48603924639782

$ws.Range("D6").Value = 0.5052686142478391
$ws.Range("E6").Value = 0.5052686142478391

$ws.Range("D7").Value = 0.1039466108553287
$ws.Range("E7").Value = 0.8960533891446714

$ws.Range("D10").Value = 0.3812488154203604
$ws.Range("E10").Value = 0.6187511845796396

$ws.Range("F11").Value = 181.4903106689453
